$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The evaluation window shifted by one additional observation: each row's
# B:G values move up into the previous row (row N <- row N+1), and the
# final row (11) receives the newly computed values for the extended
# ifoCAST series.

for ($r = 2; $r -le 10; $r++) {
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $ws.Cells.Item($r + 1, $c).Value2
    }
}

$newRow11 = @(0.1087828097550494, 0.3930655897532851, 0.2412922028369822, 0.4912150270879161, 0.5080756007043895, 9)
for ($c = 2; $c -le 7; $c++) {
    $ws.Cells.Item(11, $c).Value2 = $newRow11[$c - 2]
}
